$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 412 ("「何事も難しい。簡単になる前は」...") which shifts
# all subsequent rows up by one (row 413 becomes 412, ..., row 616 becomes 615).
$ws.Rows.Item(412).Delete()
